$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# Rows 13-27: lowercase the short "name" codes in column B.
$ws.Range("B13").Value = "gj"
$ws.Range("B14").Value = "zk"
$ws.Range("B15").Value = "ze"
$ws.Range("B16").Value = "zf"
$ws.Range("B17").Value = "za"
$ws.Range("B18").Value = "zb"
$ws.Range("B19").Value = "fs"
$ws.Range("B20").Value = "fu"
$ws.Range("B21").Value = "fp"
$ws.Range("B22").Value = "kd"
$ws.Range("B23").Value = "km"
$ws.Range("B24").Value = "kmt"
$ws.Range("B25").Value = "kmf"
$ws.Range("B26").Value = "mna"
$ws.Range("B27").Value = "mk"

# Rows 19-21: capitalize the first letter of the label text in column C.
$ws.Range("C19").Value = "Gesättigte Fettsäuren [g/Tag]"
$ws.Range("C20").Value = "Einfach ungesättigte Fettsaeuren [g/Tag]"
$ws.Range("C21").Value = "Mehrfach ungesättigte Fettsaeuren [g/Tag]"

# Update the active selection shown when the file was last saved.
$ws.Range("G27").Select()

$wb.Save()
